$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row index, Actual Consumption (MW), Timestamp (date serial)
$data = @(
    @(2, 5052, 45747.04166666666),
    @(3, 5040, 45747.05208333334),
    @(4, 5057, 45747.0625),
    @(5, 5000, 45747.07291666666),
    @(6, 5025, 45747.08333333334),
    @(7, 5005, 45747.09375),
    @(8, 5005, 45747.10416666666),
    @(9, 4977, 45747.11458333334),
    @(10, 5034, 45747.125),
    @(11, 5088, 45747.13541666666),
    @(12, 5087, 45747.14583333334),
    @(13, 5131, 45747.15625),
    @(14, 5219, 45747.16666666666),
    @(15, 5301, 45747.17708333334),
    @(16, 5345, 45747.1875),
    @(17, 5476, 45747.19791666666),
    @(18, 5751, 45747.20833333334),
    @(19, 5900, 45747.21875),
    @(20, 6012, 45747.22916666666),
    @(21, 6232, 45747.23958333334),
    @(22, 6635, 45747.25),
    @(23, 6777, 45747.26041666666),
    @(24, 6824, 45747.27083333334),
    @(25, 6958, 45747.28125),
    @(26, 7146, 45747.29166666666),
    @(27, 7282, 45747.30208333334),
    @(28, 7262, 45747.3125),
    @(29, 7321, 45747.32291666666),
    @(30, 7280, 45747.33333333334),
    @(31, 7247, 45747.34375),
    @(32, 7132, 45747.35416666666),
    @(33, 7114, 45747.36458333334),
    @(34, 7014, 45747.375),
    @(35, 6966, 45747.38541666666),
    @(36, 6814, 45747.39583333334),
    @(37, 6790, 45747.40625),
    @(38, 6619, 45747.41666666666)
)

foreach ($item in $data) {
    $r = $item[0]
    $val = $item[1]
    $ts = $item[2]
    $ws.Cells.Item($r, 1).Value = $val
    $ws.Cells.Item($r, 2).Value = $ts
}

# Remove the now-unused rows 39-56 so the sheet dimension shrinks to A1:B38
$rowsToDelete = $ws.Range("A39:A56").EntireRow
$rowsToDelete.Delete()
